$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 9

function Set-TextCell($ws, $r, $c, $text) {
    $cell = $ws.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Cells.Item($row, 1).Value = 131157738
$ws.Cells.Item($row, 2).Value = 83223
$ws.Cells.Item($row, 4).Value = "NT"
$ws.Cells.Item($row, 5).Value = 6440
$ws.Cells.Item($row, 6).Value = "Vitgrynig nållav"
$ws.Cells.Item($row, 7).Value = "Chaenotheca subroscida"
$ws.Cells.Item($row, 8).Value = "(Eitner) Zahlbr."
$ws.Cells.Item($row, 16).Value = "Edetjärnens badplats, Edetjärnens badplats, Jmt"
$ws.Cells.Item($row, 17).Value = 447575
$ws.Cells.Item($row, 18).Value = 7025519
$ws.Cells.Item($row, 19).Value = 3
$ws.Cells.Item($row, 20).Value = "Jämtland"
$ws.Cells.Item($row, 21).Value = "Krokom"
$ws.Cells.Item($row, 22).Value = "Jämtland"
$ws.Cells.Item($row, 23).Value = "Alsen"
# Startdatum/Slutdatum look like ISO dates, so Excel would normally convert
# them to date serials on assignment - force them to stay plain text instead.
Set-TextCell $ws $row 25 "2026-02-15"
$ws.Cells.Item($row, 26).Value = "11:42"
Set-TextCell $ws $row 27 "2026-02-15"
$ws.Cells.Item($row, 28).Value = "11:42"
$ws.Cells.Item($row, 29).Value = "På gammal gran"
$ws.Cells.Item($row, 30).Value = $false
$ws.Cells.Item($row, 31).Value = $false
$ws.Cells.Item($row, 33).Value = $false
$ws.Cells.Item($row, 49).Value = "Ludvig Nordin"
$ws.Cells.Item($row, 50).Value = "Ludvig Nordin"
